$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51; existing rows 51-64 shift down to 52-65.
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new weekly price record.
$ws.Cells.Item(51, 1).Value = 1
$ws.Cells.Item(51, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(51, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(51, 4).Value = 44642
$ws.Cells.Item(51, 5).Value = 15
$ws.Cells.Item(51, 6).Value = "Fruta"
$ws.Cells.Item(51, 7).Value = 100103
$ws.Cells.Item(51, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(51, 9).Value = 100103006
$ws.Cells.Item(51, 10).Value = "Nectarín"
$ws.Cells.Item(51, 11).Value = "Artic Snow"
$ws.Cells.Item(51, 12).Value = "Segunda"
$ws.Cells.Item(51, 13).Value = 300
$ws.Cells.Item(51, 14).Value = 18000
$ws.Cells.Item(51, 15).Value = 20000
$ws.Cells.Item(51, 16).Value = 19000
$ws.Cells.Item(51, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(51, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(51, 19).Value = 1056
$ws.Cells.Item(51, 20).Value = 18
